# PWR_Board_TestReportTemplate2.xlsx edit script
# "Fixed instrument connection refresh on GUI selftester temp check fixed.
#  selftester other diodes implemented."

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Quantities sheet: just a selection/view-state change (E8 -> C7)
# Do this FIRST, then come back to HWCheck so HWCheck ends up as the
# active/tabSelected sheet (matches the unchanged activeTab in workbook.xml).
# ---------------------------------------------------------------------------
$wsQuantities = $wb.Worksheets.Item("Quantities")
$wsQuantities.Range("C7").Select()

# ---------------------------------------------------------------------------
# HWCheck sheet: the bulk of the data/formula changes.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("HWCheck")
$ws.Activate()

# --- Row 22: single-cell tolerance formulas 0.95 -> 0.9 / 1.05 -> 1.1 ---
$ws.Range("D22").Formula = "=0.9*F22"
$ws.Range("E22").Formula = "=1.1*F22"

# --- Rows 23:26: shared-range shrinks from D23:D36/E23:E36 down to
#     D23:D26/E23:E26, coefficient 0.95 -> 0.9 and 1.05 -> 1.1 ---
$ws.Range("D23:D26").Formula = "=0.9*F23"
$ws.Range("E23:E26").Formula = "=1.1*F23"

# --- Rows 32:36: become their OWN shared-formula group, keeping the
#     ORIGINAL coefficients (0.95 / 1.05) - values are unchanged ---
$ws.Range("D32:D36").Formula = "=0.95*F32"
$ws.Range("E32:E36").Formula = "=1.05*F32"

# --- Row 37: tolerance formula 0.95 -> 0.9 / 1.05 -> 1.1 ---
$ws.Range("D37").Formula = "=0.9*F37"
$ws.Range("E37").Formula = "=1.1*F37"

# --- Rows 38:41: shared range keeps its extent, coefficient 0.95 -> 0.9
#     and 1.05 -> 1.1 ---
$ws.Range("D38:D41").Formula = "=0.9*F38"
$ws.Range("E38:E41").Formula = "=1.1*F38"

# --- Row 48: tolerance bumped from 1E-3 to 2E-3 ---
$ws.Range("E48").Value = 0.002

# --- Row 51: limits/measurement reworked ---
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 50
$ws.Range("F51").Value = 5

# --- Row 52: was a fixed 0 / 1E-3 pass-band, now a computed 0.9x/1.1x
#     tolerance like the rows above (style switches from the "General"
#     center style to the "0.000" center style used by D2/D47/etc). ---
$ws.Range("D52").NumberFormat = "0.000"
$ws.Range("E52").NumberFormat = "0.000"
$ws.Range("D52").Formula = "=0.9*F52"
$ws.Range("E52").Formula = "=1.1*F52"
$ws.Range("F52").Value = 4.5

# --- Rows 53:56: new "P/F" (Pass/Fail) column label in H, diode selftest ---
$ws.Range("H53").HorizontalAlignment = -4108
$ws.Range("H54").HorizontalAlignment = -4108
$ws.Range("H55").HorizontalAlignment = -4108
$ws.Range("H56").HorizontalAlignment = -4108
$ws.Range("H53").Value = "P/F"
$ws.Range("H54").Value = "P/F"
$ws.Range("H55").Value = "P/F"
$ws.Range("H56").Value = "P/F"

# --- Restore the view/selection state on HWCheck (frozen header pane,
#     scrolled down, G42 -> F51 active cell) ---
$win = $excel.ActiveWindow
$win.FreezePanes = $true
$ws.Range("F51").Select()
